{"js": "// Code-review checklist header table update:\n//   Sprint No.  : \"1\"        -> \"2\"\n//   Review Date : \"02/09/18\" -> \"02/21/18\" (appears in the merged date cell)\n//\n// The table layout (row index, cell index) is:\n//   row 0: Project Name | <value>      | Project ID    | <value>\n//   row 1: Reviewer's Name | <value>   | Sprint No.    | <value>  <- target\n//   row 2: Review Date | <value, gridSpan=3>                      <- target\n//   row 3: File Name (Source Code) | <value> | <value> | <value>\n//\n// We scope the text search to the specific target cell so that other runs\n// containing the same literal text elsewhere in the document (e.g. the\n// \"1.  License\" heading) are left untouched, and we replace only the\n// matched range (not the whole cell) so existing run/paragraph formatting\n// is preserved.\n\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// --- Sprint No.: \"1\" -> \"2\" ----------------------------------------------\nconst sprintRow = rows.items[1];\nconst sprintCells = sprintRow.cells;\nsprintCells.load(\"items\");\nawait context.sync();\n\nconst sprintValueCell = sprintCells.items[3];\nconst sprintMatches = sprintValueCell.body.search(\"1\", {\n  matchCase: true,\n  matchWholeWord: true\n});\nsprintMatches.load(\"items\");\nawait context.sync();\n\nif (sprintMatches.items.length > 0) {\n  sprintMatches.items[0].insertText(\"2\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- Review Date: \"02/09/18\" -> \"02/21/18\" --------------------------------\nconst dateRow = rows.items[2];\nconst dateCells = dateRow.cells;\ndateCells.load(\"items\");\nawait context.sync();\n\nconst dateValueCell = dateCells.items[1];\nconst dateMatches = dateValueCell.body.search(\"02/09/18\", {\n  matchCase: true,\n  matchWholeWord: false\n});\ndateMatches.load(\"items\");\nawait context.sync();\n\nif (dateMatches.items.length > 0) {\n  dateMatches.items[0].insertText(\"02/21/18\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Code-review checklist header table update:\n#   Sprint No.  : \"1\"        -> \"2\"\n#   Review Date : \"02/09/18\" -> \"02/21/18\"\n#\n# Table 1 (header/info table), 1-based Cell(row, col) as COM indexes it:\n#   Row 1: Project Name     | <value>                    | Project ID  | <value>\n#   Row 2: Reviewer's Name  | <value>                    | Sprint No.  | <value>  <- Cell(2,4)\n#   Row 3: Review Date      | <value, spans columns 2-4>                          <- Cell(3,2)\n#   Row 4: File Name (Source Code) | <value> | <value>   | <value>\n#\n# We scope Find to each specific target cell's Range so that other runs with\n# the same literal text elsewhere in the document (e.g. the \"1.  License\"\n# heading) are left untouched. After Find locates the match, the text is\n# replaced by assigning Range.Text directly (rather than via Find's\n# Replace argument) so the existing run formatting / xml:space is preserved\n# exactly.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# --- Sprint No.: \"1\" -> \"2\" ------------------------------------------------\n$sprintCell = $table.Cell(2, 4)\n$sprintRange = $sprintCell.Range\n$sprintFind = $sprintRange.Find\n$sprintFind.ClearFormatting()\n$sprintFind.Execute(\"1\", $true, $true, $false, $false, $false, $true, 1, $false) | Out-Null\n$sprintRange.Text = \"2\"\n\n# --- Review Date: \"02/09/18\" -> \"02/21/18\" ---------------------------------\n$dateCell = $table.Cell(3, 2)\n$dateRange = $dateCell.Range\n$dateFind = $dateRange.Find\n$dateFind.ClearFormatting()\n$dateFind.Execute(\"02/09/18\", $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null\n$dateRange.Text = \"02/21/18\"\n"}
